# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) for the last row (row 18)
# of both the "zh-cn" and "de-de" report sheets, reflecting a
# freshly regenerated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet, row 18: D = Correspond Handoff Datetime, G = Correspond Handback DateTime
$wsZhCn.Cells.Item(18, 4).Value = "2016-03-01 09:44:23"
$wsZhCn.Cells.Item(18, 7).Value = "2016-03-01 09:45:09"

# de-de sheet, row 18: D = Correspond Handoff Datetime, G = Correspond Handback DateTime
$wsDeDe.Cells.Item(18, 4).Value = "2016-03-01 09:44:34"
$wsDeDe.Cells.Item(18, 7).Value = "2016-03-01 09:45:28"
